$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 19/20: Chainlink and Polkadot swap positions, with new price/volume values
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'5.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.79%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'13.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.81%  "
$ws.Range("E20").Style = "Normal"

# Price / Volume updates for remaining rows
$ws.Range("D2").Value = "'61.426.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.38%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.321.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.87%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'566.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.64%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'128.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.08%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.321.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.82%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.78%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -3.86%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.377"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.57%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.887.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.79%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.29%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.324.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.71%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -5.15%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'24.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.05%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'61.480.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.19%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D21").Value = "'9.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -9.74%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'356.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -7.26%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -3.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.452.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.90%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'69.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.52%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -5.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.64%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.62%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.61%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -5.39%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.54%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.349.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.86%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'22.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.23%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'161.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0762"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'4.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'41.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -7.47%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.68%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -4.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'6.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'22.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -8.23%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.862"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -6.53%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'21.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.63%  "
$ws.Range("E51").Style = "Normal"
